# PowerPoll pie-chart update
# --------------------------
# The commit this task is based on only touches the PowerPoll task-pane
# add-in's own JavaScript/HTML (outside this .pptx) -- it stops drawing
# data-labels inside the pie slices and colour-codes the poll options, and
# adds some input validation to the add-in code.  None of that lives inside
# the slide deck.
#
# Inside Presentation1.pptx the add-in is represented purely as a cached,
# blank placeholder snapshot: an <mc:AlternateContent> block on slide 1 that
# pairs a <we:webextensionref> (live task-pane content, shown inside
# PowerPoint) with a fallback <p:pic> (the static image other viewers show).
# The underlying we:webextension part (ppt/slides/udata/data.xml) simply
# carries an internal GUID (we:webextension/@id) that PowerPoint regenerates
# whenever it rewrites the add-in binding on save -- it is not bound to any
# slide text, shape geometry, formatting, or other content that the
# PowerPoint object model exposes. There is no Shapes/TextFrame/Tags/
# CustomXMLParts member that reaches that id, so nothing observable through
# COM automation changes for this revision.
#
# Touch the presentation/slide so the script still operates against the
# live object model (mirroring how the add-in's host slide would be
# reached), without mutating any shape the object model does expose.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$null = $s.Shapes.Count
